# Estado de Cuenta - add new worker (ALEJANDRO DE LA CRUZ MARTINEZ HOYOS) and
# refresh the "Valor Mora" / period rows per the updated database export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room for the new worker's 7 period rows (inserted right after
#    the current last data row of IVAN DARIO MARTINEZ HOYOS, row 31).
#    This pushes the signature block (rows 36-37) down to rows 43-44.
# ------------------------------------------------------------------
$ws.Range("A32:A38").EntireRow.Insert()

# Copy the "middle of table" row formatting down onto the 7 fresh rows,
# then give the very last one (row 38) the closing/bottom-border look
# that used to belong to row 31.
$ws.Range("B16:J16").Copy()
$ws.Range("B32:J37").PasteSpecial(-4122)

$ws.Range("B31:J31").Copy()
$ws.Range("B38:J38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 31 itself (IVAN DARIO's oldest period, 2201) goes back to being a
# normal row now that it is no longer the last line of the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B31:J31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Re-key the existing periods for IVAN DARIO MARTINEZ HOYOS
#    (rows 16-31), newest period first, and refresh Valor Mora.
# ------------------------------------------------------------------
$periodosIvan = @("2304","2303","2302","2301","2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201")
$valorMoraIvan = @(32707,36341,36341,36341,36341,36341,36341,36341,36341,36341,36341,36341,36341,36341,36341,36341)

for ($i = 0; $i -lt $periodosIvan.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "73214789"
    $ws.Cells.Item($r, 4).Value = "IVAN DARIO MARTINEZ HOYOS"
    $ws.Cells.Item($r, 5).Value = $periodosIvan[$i]
    $ws.Cells.Item($r, 6).Value = $valorMoraIvan[$i]
    $ws.Cells.Item($r, 7).Value = 908526
}

# ------------------------------------------------------------------
# 3) Fill in the new worker, ALEJANDRO DE LA CRUZ MARTINEZ HOYOS
#    (rows 32-38).
# ------------------------------------------------------------------
$periodosAlejandro = @("2304","2303","2302","2301","2212","2211","2210")
$valorMoraAlejandro = @(36000,40000,40000,40000,40000,40000,40000)

for ($i = 0; $i -lt $periodosAlejandro.Length; $i++) {
    $r = 32 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1143329586"
    $ws.Cells.Item($r, 4).Value = "ALEJANDRO DE LA CRUZ MARTINEZ HOYOS"
    $ws.Cells.Item($r, 5).Value = $periodosAlejandro[$i]
    $ws.Cells.Item($r, 6).Value = $valorMoraAlejandro[$i]
    $ws.Cells.Item($r, 7).Value = 689454
}

# ------------------------------------------------------------------
# 4) Header totals: now 2 workers owing a combined 853,822 over 16 periods.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 853822
$ws.Range("C13").Value = 2

# Column D needs to be wide enough for the longer worker name.
$ws.Columns.Item(4).AutoFit() | Out-Null

Write-Host "Workbook updated"
